$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-101) forward by 3 days
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 3
}

# Update Excedent_RO_UA (J), Total Excedent (L), System Direction (N) for rows 50-53
for ($r = 50; $r -le 53; $r++) {
    $ws.Cells.Item($r, 10).Value2 = 4.25
    $ws.Cells.Item($r, 12).Value2 = 4.25
    $ws.Cells.Item($r, 14).Value2 = 4.25
}

# Update same columns for rows 54-57
for ($r = 54; $r -le 57; $r++) {
    $ws.Cells.Item($r, 10).Value2 = 173.25
    $ws.Cells.Item($r, 12).Value2 = 173.25
    $ws.Cells.Item($r, 14).Value2 = 173.25
}
